# Automatic daily update of the "remaining days" (剩余, column E) tracker.
# For every data row (starting at row 2):
#   - If the remaining-days counter (E) is greater than 1, decrement it by 1.
#   - If the remaining-days counter (E) equals 1, the cycle has completed:
#       reset E back to the total-days value (D) and roll the start date (F)
#       forward to the new cycle start date (20260127).
#   - Rows whose start date (F) is not a well-formed 8-digit yyyymmdd value
#     (data-entry error) are skipped and left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCycleDate = 20260127

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D - total days
    $eCell = $ws.Cells.Item($r, 5)   # column E - remaining days
    $fCell = $ws.Cells.Item($r, 6)   # column F - start date

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null) {
        continue
    }

    # Skip rows with a malformed start date (not an 8-digit yyyymmdd number).
    if ($fVal -eq $null -or "$fVal".Length -ne 8) {
        continue
    }

    if ($eVal -le 1) {
        $eCell.Value2 = $dVal
        $fCell.Value2 = $newCycleDate
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
